$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4247.0586
$ws.Range("I76").Value = 3652.3809
$ws.Range("J76").Value = 5207.6924
$ws.Range("K76").Value = 3652.3809
$ws.Range("L76").Value = 5207.6924
$ws.Range("M76").Value = -3337.3809
$ws.Range("N76").Value = -5837.6924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4247.0586
$ws.Range("I79").Value = 3652.3809
$ws.Range("J79").Value = 5207.6924
$ws.Range("K79").Value = 3652.3809
$ws.Range("L79").Value = 5207.6924
$ws.Range("M79").Value = -2560.3809
$ws.Range("N79").Value = -7391.6924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 11013.786
$ws.Range("I88").Value = 2932.3333
$ws.Range("J88").Value = 13217.818
$ws.Range("K88").Value = 2932.3333
$ws.Range("L88").Value = 13217.818
$ws.Range("M88").Value = -2526.3333
$ws.Range("N88").Value = -14029.818

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 11013.786
$ws.Range("I91").Value = 2932.3333
$ws.Range("J91").Value = 13217.818
$ws.Range("K91").Value = 2932.3333
$ws.Range("L91").Value = 13217.818
$ws.Range("M91").Value = -1528.3333
$ws.Range("N91").Value = -16025.818

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5052220.5
$ws.Range("I132").Value = 1265.9656
$ws.Range("J132").Value = 41671640
$ws.Range("K132").Value = 3797.8968
$ws.Range("L132").Value = 125014920
$ws.Range("M132").Value = -1267.8968
$ws.Range("N132").Value = -125019980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6719.785
$ws.Range("I32").Value = 4973.025
$ws.Range("K32").Value = 4973.025
$ws.Range("M32").Value = -4686.025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1850.8
$ws.Range("I88").Value = 1814
$ws.Range("J88").Value = 1998
$ws.Range("K88").Value = 1814
$ws.Range("L88").Value = 1998
$ws.Range("M88").Value = -1408
$ws.Range("N88").Value = -2810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1850.8
$ws.Range("I91").Value = 1814
$ws.Range("J91").Value = 1998
$ws.Range("K91").Value = 1814
$ws.Range("L91").Value = 1998
$ws.Range("M91").Value = -410
$ws.Range("N91").Value = -4806

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1615717.5
$ws.Range("I132").Value = 1992.8372
$ws.Range("J132").Value = 5267831.5
$ws.Range("K132").Value = 5978.5116
$ws.Range("L132").Value = 15803494.5
$ws.Range("M132").Value = -3448.5116
$ws.Range("N132").Value = -15808554.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 900
$ws.Range("I7").Value = 900
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -787

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1575.7693
$ws.Range("I86").Value = 1507.7273
$ws.Range("J86").Value = 1950
$ws.Range("K86").Value = 1507.7273
$ws.Range("L86").Value = 1950
$ws.Range("M86").Value = -384.7273
$ws.Range("N86").Value = -4196

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1575.7693
$ws.Range("I89").Value = 1507.7273
$ws.Range("J89").Value = 1950
$ws.Range("K89").Value = 7538.636500000001
$ws.Range("L89").Value = 9750
$ws.Range("M89").Value = -1922.636500000001
$ws.Range("N89").Value = -20982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1000000000
$ws.Range("I23").Value = 1000000000
$ws.Range("K23").Value = 1000000000
$ws.Range("M23").Value = -999999760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 6110
$ws.Range("J26").Value = 6110
$ws.Range("L26").Value = 6110
$ws.Range("N26").Value = -6684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 1000000000
$ws.Range("I27").Value = 1000000000
$ws.Range("K27").Value = 1000000000
$ws.Range("M27").Value = -999999808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("N52").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6796618.5
$ws.Range("I58").Value = 7938241
$ws.Range("J58").Value = 2001805.8
$ws.Range("K58").Value = 7938241
$ws.Range("L58").Value = 2001805.8
$ws.Range("M58").Value = -7938038
$ws.Range("N58").Value = -2002211.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9221.857
$ws.Range("I62").Value = 9960
$ws.Range("J62").Value = 7376.5
$ws.Range("K62").Value = 9960
$ws.Range("L62").Value = 7376.5
$ws.Range("M62").Value = -9336
$ws.Range("N62").Value = -8624.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 9221.857
$ws.Range("I65").Value = 9960
$ws.Range("J65").Value = 7376.5
$ws.Range("K65").Value = 49800
$ws.Range("L65").Value = 36882.5
$ws.Range("M65").Value = -46680
$ws.Range("N65").Value = -43122.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1077.9445
$ws.Range("I107").Value = 700.2143
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 700.2143
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 1219.7857
$ws.Range("N107").Value = -6240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6796618.5
$ws.Range("I136").Value = 7938241
$ws.Range("J136").Value = 2001805.8
$ws.Range("K136").Value = 23814723
$ws.Range("L136").Value = 6005417.4
$ws.Range("M136").Value = -23812173
$ws.Range("N136").Value = -6010517.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4984.5293
$ws.Range("I70").Value = 4775.7
$ws.Range("J70").Value = 5282.857
$ws.Range("K70").Value = 4775.7
$ws.Range("L70").Value = 5282.857
$ws.Range("M70").Value = -4505.7
$ws.Range("N70").Value = -5822.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4984.5293
$ws.Range("I73").Value = 4775.7
$ws.Range("J73").Value = 5282.857
$ws.Range("K73").Value = 4775.7
$ws.Range("L73").Value = 5282.857
$ws.Range("M73").Value = -3839.7
$ws.Range("N73").Value = -7154.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8887
$ws.Range("I80").Value = 16157.857
$ws.Range("J80").Value = 2525
$ws.Range("K80").Value = 16157.857
$ws.Range("L80").Value = 2525
$ws.Range("M80").Value = -15159.857
$ws.Range("N80").Value = -4521

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 8887
$ws.Range("I83").Value = 16157.857
$ws.Range("J83").Value = 2525
$ws.Range("K83").Value = 80789.285
$ws.Range("L83").Value = 12625
$ws.Range("M83").Value = -75797.285
$ws.Range("N83").Value = -22609

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2243.7144
$ws.Range("I102").Value = 1706.4286
$ws.Range("J102").Value = 3318.2856
$ws.Range("K102").Value = 1706.4286
$ws.Range("L102").Value = 3318.2856
$ws.Range("M102").Value = -84.42859999999996
$ws.Range("N102").Value = -6562.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 47620108
$ws.Range("I113").Value = 90909890
$ws.Range("J113").Value = 1349
$ws.Range("K113").Value = 90909890
$ws.Range("L113").Value = 1349
$ws.Range("M113").Value = -90907720
$ws.Range("N113").Value = -5689

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 71834380
$ws.Range("I122").Value = 81910620
$ws.Range("J122").Value = 50002532
$ws.Range("K122").Value = 245731860
$ws.Range("L122").Value = 150007596
$ws.Range("M122").Value = -245729410
$ws.Range("N122").Value = -150012496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 13123.235
$ws.Range("I126").Value = 16084.615
$ws.Range("J126").Value = 3498.75
$ws.Range("K126").Value = 48253.845
$ws.Range("L126").Value = 10496.25
$ws.Range("M126").Value = -45783.845
$ws.Range("N126").Value = -15436.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 56666.668
$ws.Range("J136").Value = 56666.668
$ws.Range("L136").Value = 170000.004
$ws.Range("N136").Value = -175100.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 855354.25
$ws.Range("I82").Value = 2001500.8
$ws.Range("J82").Value = 139012.62
$ws.Range("K82").Value = 2001500.8
$ws.Range("L82").Value = 139012.62
$ws.Range("M82").Value = -2001139.8
$ws.Range("N82").Value = -139734.62

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 855354.25
$ws.Range("I85").Value = 2001500.8
$ws.Range("J85").Value = 139012.62
$ws.Range("K85").Value = 2001500.8
$ws.Range("L85").Value = 139012.62
$ws.Range("M85").Value = -2000252.8
$ws.Range("N85").Value = -141508.62

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4291209
$ws.Range("I122").Value = 5500943.5
$ws.Range("J122").Value = 1670117.5
$ws.Range("K122").Value = 16502830.5
$ws.Range("L122").Value = 5010352.5
$ws.Range("M122").Value = -16500380.5
$ws.Range("N122").Value = -5015252.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 36360.875
$ws.Range("J46").Value = 36360.875
$ws.Range("L46").Value = 36360.875
$ws.Range("N46").Value = -36822.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16751467
$ws.Range("I81").Value = 1400.1428
$ws.Range("J81").Value = 40201560
$ws.Range("K81").Value = 2800.2856
$ws.Range("L81").Value = 80403120
$ws.Range("M81").Value = -1739.2856
$ws.Range("N81").Value = -80405242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 16751467
$ws.Range("I84").Value = 1400.1428
$ws.Range("J84").Value = 40201560
$ws.Range("K84").Value = 14001.428
$ws.Range("L84").Value = 402015600
$ws.Range("M84").Value = -8697.428
$ws.Range("N84").Value = -402026208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 36360.875
$ws.Range("J134").Value = 36360.875
$ws.Range("L134").Value = 109082.625
$ws.Range("N134").Value = -114152.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4765602.5
$ws.Range("I136").Value = 4581.6875
$ws.Range("J136").Value = 8774883
$ws.Range("K136").Value = 13745.0625
$ws.Range("L136").Value = 26324649
$ws.Range("M136").Value = -11195.0625
$ws.Range("N136").Value = -26329749
